$wb = $excel.ActiveWorkbook

# "0baterias" sheet holds the initial lithium-ion battery count; update it to 232
$ws = $wb.Worksheets.Item("0baterias")
$ws.Range("B2").Value = 232

# Reflect the selection recorded in the saved file (active cell G12 on this sheet)
$ws.Activate()
$ws.Range("G12").Select()
